$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-10 from 2023-10-25 (45224) to 2023-11-03 (45233)
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45233
}
